$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "26.659.48", "1.001") are preserved exactly as text, matching the
# original inlineStr cell content instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.659.48"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.792.88"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "309.12"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.4447"
$ws.Range("E7").Value = "  +5.06%  "
$ws.Range("D8").Value = "0.3656"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "0.07284"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "0.8529"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").Value = "20.55"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "1.923.28"
$ws.Range("E12").Value = "  +5.34%  "
$ws.Range("D13").Value = "6.602"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "0.07068"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "91.86"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "5.256"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "0.000008651"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "14.76"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "26.697.71"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").Value = "5.132"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "10.75"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "151.79"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.172"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.34"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "5.165"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "116.49"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "0.08777"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "0.7395"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "1.153"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").Value = "2.924"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "4.437"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "1.084"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").Value = "0.01953"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "0.05154"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").Value = "0.5272"
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("D40").Value = "2.831"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").Value = "7.001"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "0.5082"
$ws.Range("E43").Value = "  +7.17%  "
$ws.Range("D44").Value = "8.387"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "10.40"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").Value = "1.949"
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("D47").Value = "105.17"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "1.651"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "0.06290"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "0.9100"
$ws.Range("E51").Value = "  -0.05%  "

# Restore column D to the default "Normal" style/format so no stray cell
# formatting is introduced (values remain text because they are already
# stored as string cells).
$ws.Range("D2:D51").Style = "Normal"
